# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-04 08:22:36
#
# Applies the data corrections to the "Session Analysis Results" sheet:
#  - normalizes a handful of "Recorded By" email lists (re-ordering only)
#  - updates Class / Group statistics figures (K:S columns)
#  - flips four sessions (rows 46, 107, 142, 164) from Pending/Not Recorded
#    to Recorded, filling in their "Recorded By" / "Students" / "Status"
#    cells and matching the green "Recorded" row styling
#  - refreshes a few "Students" attendance counts (H column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple value-only updates: [cellRef, newValue]
# ---------------------------------------------------------------------
$valueUpdates = @(
    @("G2",  "majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"),

    @("L6",  36),
    @("L7",  11),
    @("L8",  129),
    @("L9",  "20.5%"),
    @("L10", "30.9%"),

    @("S16", "37.2%"),

    @("O17", 5),
    @("P17", 1),
    @("Q17", 16),
    @("R17", "22.7%"),
    @("S17", "37.3%"),

    @("G18", "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"),

    @("G19", "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"),
    @("O19", 4),
    @("P19", 3),
    @("Q19", 15),
    @("R19", "18.2%"),
    @("S19", "30.5%"),

    @("O21", 5),
    @("P21", 0),
    @("Q21", 17),
    @("R21", "22.7%"),
    @("S21", "30.6%"),

    @("O22", 5),
    @("P22", 1),
    @("Q22", 16),
    @("R22", "22.7%"),
    @("S22", "13.5%"),

    @("G24", "majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"),
    @("H24", "153/217"),

    @("G40", "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"),

    @("G41", "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"),

    @("G54", "merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"),

    @("G58", "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"),

    @("G62", "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"),

    @("G76", "merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"),

    @("G80", "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"),

    @("G84", "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"),

    @("G98", "merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"),

    @("G106", "neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"),

    @("G120", "merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"),

    @("G128", "neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"),

    @("G134", "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"),
    @("H134", "84/224"),

    @("G150", "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"),
    @("H150", "95/224"),

    @("G156", "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"),

    @("G172", "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg")
)

foreach ($upd in $valueUpdates) {
    $ws.Range($upd[0]).Value = $upd[1]
}

# ---------------------------------------------------------------------
# 2) Rows that flip from Pending / Not Recorded -> Recorded.
#    Copy the green "Recorded" formatting (row 8 is a clean s="2" donor
#    row spanning A:I) onto each target row, then write the new
#    Recorded-By / Students / Status text. Columns A-F keep their
#    existing text, only their fill/style changes.
# ---------------------------------------------------------------------
$ws.Range("A8:I8").Copy() | Out-Null
$recordedRows = @(46, 107, 142, 164)
foreach ($r in $recordedRows) {
    $ws.Range("A$r`:I$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

$rowUpdates = @(
    @(46,  "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg", "51/220", "Recorded"),
    @(107, "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg", "33/154", "Recorded"),
    @(142, "merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg", "100/224", "Recorded"),
    @(164, "merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg", "1/226", "Recorded")
)

foreach ($ru in $rowUpdates) {
    $r = $ru[0]
    $ws.Cells.Item($r, 7).Value = $ru[1]   # G - Recorded By
    $ws.Cells.Item($r, 8).Value = $ru[2]   # H - Students
    $ws.Cells.Item($r, 9).Value = $ru[3]   # I - Status
}
